# Refresh the "Saldo" export:
#  - drop the LAIS row (account closed / no longer in the filtered export)
#  - drop the DAIANNE row, and bump GUILHERME's balance (877.48 -> 1524.94),
#    re-inserting that account where its new, higher balance sorts it
#  - bump BLUEMETRIX's balance (946.85 -> 953.9), re-inserting it just above
#    ASPA where its new balance sorts it
#
# Work from the bottom of the sheet upward so that row numbers referenced
# below remain valid as earlier rows are deleted/inserted.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22: old GUILHERME (004574428 / 877.48) is superseded by the updated
# row we will write at row 9 below -- remove the stale one.
$ws.Rows.Item(22).Delete()

# Row 13: old BLUEMETRIX (001761119 / 946.85) is superseded by the updated
# row we will insert above ASPA -- remove the stale one.
$ws.Rows.Item(13).Delete()

# Insert the refreshed BLUEMETRIX row just above ASPA (row 12) with its new
# balance of 953.9.
$ws.Rows.Item(12).Insert()
$ws.Cells.Item(12, 1).Value = "'001761119"
$ws.Cells.Item(12, 2).Value = "BLUEMETRIX"
$ws.Cells.Item(12, 3).Value = 953.9

# Row 9 previously held DAIANNE (004473942 / 3839.36); overwrite it in place
# with GUILHERME's refreshed balance (004574428 / 1524.94).
$ws.Cells.Item(9, 1).Value = "'004574428"
$ws.Cells.Item(9, 2).Value = "GUILHERME"
$ws.Cells.Item(9, 3).Value = 1524.94

# Row 3: LAIS (004230529 / 45901.8) no longer appears in the export -- remove it.
$ws.Rows.Item(3).Delete()
